# Add a "Category" header in A1 (matching the style of the other header
# cells in row 1), and strip the bold/header formatting that was
# previously (mistakenly) applied to the category cells in A2:A46 so
# that only the new header row carries that style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header-style formatting from the category column's data rows
# (A2:A46) - they become plain/default-styled cells.
$ws.Range("A2:A46").ClearFormats()

# Insert the new column header text.
$ws.Range("A1").Value = "Category"

# Give A1 the same formatting as the rest of the header row (bold,
# centered, top-aligned, bordered) by copying the format from B1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
